# Last adjustments before v2 - update recomputed patient metrics
# (points selected in the ECG curve + downstream D-station parameters)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("U3").Value = 160
$ws.Range("V3").Value = 1016
$ws.Range("W3").Value = 1195
$ws.Range("X3").Value = 160
$ws.Range("AC3").Value = 1016
$ws.Range("AD3").Value = 1195
$ws.Range("U4").Value = 149
$ws.Range("V4").Value = 993
$ws.Range("X4").Value = 149
$ws.Range("AC4").Value = 993
$ws.Range("AJ4").Value = -7.5
$ws.Range("AK4").Value = 151.8
$ws.Range("U5").Value = 159
$ws.Range("V5").Value = 966
$ws.Range("W5").Value = 1181
$ws.Range("X5").Value = 159
$ws.Range("AC5").Value = 966
$ws.Range("AD5").Value = 1181
$ws.Range("AH5").Value = 806.9999999999999
$ws.Range("AI5").Value = 0.5018587360594796
$ws.Range("AJ5").Value = -11.6
$ws.Range("AK5").Value = 77.40000000000001
$ws.Range("U8").Value = 152
$ws.Range("V8").Value = 995
$ws.Range("W8").Value = 1120
$ws.Range("X8").Value = 152
$ws.Range("AC8").Value = 995
$ws.Range("AD8").Value = 1120
$ws.Range("AG8").Value = -439.0000000000001
$ws.Range("AI8").Value = -0.2787301587301588
$ws.Range("U9").Value = 127
$ws.Range("V9").Value = 820
$ws.Range("W9").Value = 1001
$ws.Range("X9").Value = 127
$ws.Range("AC9").Value = 820
$ws.Range("AD9").Value = 1001
$ws.Range("AG9").Value = -490.0000000000001
$ws.Range("AI9").Value = -0.3230059327620304
$ws.Range("U10").Value = 96
$ws.Range("V10").Value = 688
$ws.Range("W10").Value = 790
$ws.Range("X10").Value = 96
$ws.Range("AC10").Value = 688
$ws.Range("AD10").Value = 790
$ws.Range("AG10").Value = -384.9999999999999
$ws.Range("AI10").Value = -0.3202995008319467
$ws.Range("AJ10").Value = -4.6
$ws.Range("AK10").Value = 96.2
$ws.Range("U11").Value = 101
$ws.Range("V11").Value = 703
$ws.Range("W11").Value = 784
$ws.Range("X11").Value = 101
$ws.Range("AC11").Value = 703
$ws.Range("AD11").Value = 784
$ws.Range("AG11").Value = -359.0000000000001
$ws.Range("AI11").Value = -0.3068376068376069
$ws.Range("U12").Value = 144
$ws.Range("V12").Value = 951
$ws.Range("W12").Value = 1106
$ws.Range("X12").Value = 144
$ws.Range("AC12").Value = 951
$ws.Range("AD12").Value = 1106
$ws.Range("AG12").Value = 342.0000000000001
$ws.Range("AH12").Value = 778.9999999999999
$ws.Range("AI12").Value = 0.4390243902439026
$ws.Range("U13").Value = 108
$ws.Range("V13").Value = 771
$ws.Range("W13").Value = 885
$ws.Range("X13").Value = 108
$ws.Range("AC13").Value = 771
$ws.Range("AD13").Value = 885
$ws.Range("AI13").Value = -0.3423753665689149
$ws.Range("AJ13").Value = -11.6
$ws.Range("AK13").Value = 51.2
$ws.Range("U15").Value = 173
$ws.Range("V15").Value = 1167
$ws.Range("W15").Value = 1255
$ws.Range("X15").Value = 173
$ws.Range("AC15").Value = 1167
$ws.Range("AD15").Value = 1255
$ws.Range("AI15").Value = 0.4414414414414414
$ws.Range("U16").Value = 130
$ws.Range("V16").Value = 771
$ws.Range("W16").Value = 937
$ws.Range("X16").Value = 130
$ws.Range("AC16").Value = 771
$ws.Range("AD16").Value = 937
$ws.Range("AI16").Value = 0.6220338983050848
$ws.Range("U17").Value = 132
$ws.Range("V17").Value = 828
$ws.Range("W17").Value = 1022
$ws.Range("X17").Value = 132
$ws.Range("AC17").Value = 828
$ws.Range("AD17").Value = 1022
$ws.Range("AI17").Value = -0.3292993630573249
$ws.Range("U18").Value = 140
$ws.Range("V18").Value = 896
$ws.Range("W18").Value = 1001
$ws.Range("X18").Value = 140
$ws.Range("AC18").Value = 896
$ws.Range("AD18").Value = 1001
$ws.Range("AI18").Value = -0.3486005089058524
$ws.Range("U19").Value = 186
$ws.Range("V19").Value = 1076
$ws.Range("W19").Value = 1292
$ws.Range("X19").Value = 186
$ws.Range("AC19").Value = 1076
$ws.Range("AD19").Value = 1292
$ws.Range("AH19").Value = 944.0000000000001
$ws.Range("AI19").Value = 0.3983050847457627
$ws.Range("U23").Value = 162
$ws.Range("V23").Value = 1039
$ws.Range("W23").Value = 1217
$ws.Range("X23").Value = 162
$ws.Range("AC23").Value = 1039
$ws.Range("AD23").Value = 1217
$ws.Range("AG23").Value = -438.9999999999999
$ws.Range("AI23").Value = -0.2613095238095238
$ws.Range("U32").Value = 187
$ws.Range("V32").Value = 1160
$ws.Range("W32").Value = 1304
$ws.Range("X32").Value = 187
$ws.Range("AC32").Value = 1160
$ws.Range("AD32").Value = 1304
$ws.Range("AH32").Value = 950.0000000000001
$ws.Range("AI32").Value = 0.3978947368421052
$ws.Range("AJ32").Value = -16.4
$ws.Range("AK32").Value = 33.7
